$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply centered style (matches existing style used by B2:F1 header row) to the
# --- whole matrix body (B3:F9) plus the two trailing header cells F2/G2 that are
# --- new in this revision. This must happen before/after setting values; alignment
# --- is independent of cell content.
$ws.Range("B3:F9").HorizontalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("G2").HorizontalAlignment = -4108

# --- Update the class column headers (row 2) with full class names, and fill in
# --- the traceability matrix body with "X" marks, in the same order the original
# --- author appears to have performed the edits (first 3 existing class headers,
# --- then the bulk of the "X" marks, then the two new class headers that were
# --- appended at the end, then the remaining "X" marks referencing them). ---
$ws.Range("B2").Value = "CLS1: Vehiculo"
$ws.Range("C2").Value = "CLS2: Comentario"
$ws.Range("D2").Value = "CLS3: Denuncia"

$ws.Range("B3").Value = "X"
$ws.Range("B4").Value = "X"
$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("B6").Value = "X"
$ws.Range("D6").Value = "X"
$ws.Range("B7").Value = "X"
$ws.Range("B8").Value = "X"
$ws.Range("B9").Value = "X"

$ws.Range("E2").Value = "CLS4: Propietario"
$ws.Range("F2").Value = "CLS5: Infraccion"

$ws.Range("E4").Value = "X"
$ws.Range("E8").Value = "X"
$ws.Range("F8").Value = "X"
$ws.Range("E9").Value = "X"
$ws.Range("F9").Value = "X"

# --- Column widths for the newly-populated class columns ---
$ws.Columns.Item(2).ColumnWidth = 15.307291666666666
$ws.Columns.Item(3).ColumnWidth = 17.451822916666668
$ws.Columns.Item(4).ColumnWidth = 15.877604166666666
$ws.Columns.Item(5).ColumnWidth = 17.736979166666668
$ws.Columns.Item(6).ColumnWidth = 15.307291666666666

# --- Selection moved to C9 in the saved view ---
[void]$ws.Range("C9").Select()
